# Weekly data refresh: insert a new record (latest week) as row 72,
# pushing the existing rows 72:110 down to 73:111.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 72 (shifts 72:110 -> 73:111)
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new observation
$ws.Cells.Item(72, 1).Value = 10
$ws.Cells.Item(72, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(72, 3).Value = "La Araucanía"
$ws.Cells.Item(72, 4).Value = 45236
$ws.Cells.Item(72, 5).Value = 9
$ws.Cells.Item(72, 6).Value = 100112042
$ws.Cells.Item(72, 7).Value = "Locoto"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 50
$ws.Cells.Item(72, 11).Value = 3800
$ws.Cells.Item(72, 12).Value = 3800
$ws.Cells.Item(72, 13).Value = 3800
$ws.Cells.Item(72, 14).Value = "`$/kilo"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 3800
$ws.Cells.Item(72, 17).Value = 1
$ws.Cells.Item(72, 18).Value = "Hortaliza"
